# Update "想去人数" (number of people interested) figures across the
# workbook's sheets to reflect newly generated output.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 303
$ws1.Range("F5").Value  = 946
$ws1.Range("F7").Value  = 1537
$ws1.Range("F8").Value  = 39639
$ws1.Range("F11").Value = 8427
$ws1.Range("F14").Value = 699
$ws1.Range("F18").Value = 676
$ws1.Range("F22").Value = 230
$ws1.Range("F23").Value = 1073
$ws1.Range("F27").Value = 570
$ws1.Range("F34").Value = 4
$ws1.Range("F44").Value = 1042

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 4385

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 334
$ws3.Range("F5").Value = 104

# Sheet "全部类型" (All Types) - aggregate of the other sheets
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 334
$ws4.Range("F5").Value  = 303
$ws4.Range("F7").Value  = 946
$ws4.Range("F8").Value  = 1537
$ws4.Range("F9").Value  = 39639
$ws4.Range("F15").Value = 8427
$ws4.Range("F19").Value = 699
$ws4.Range("F24").Value = 676
$ws4.Range("F29").Value = 230
$ws4.Range("F30").Value = 1073
$ws4.Range("F33").Value = 570
$ws4.Range("F45").Value = 1042
